$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 695, pushing existing rows 695:732 down to 696:733
$ws.Rows.Item(695).Insert()

# Populate the newly inserted row 695 with the new record
$ws.Range("A695").Value = 11
$ws.Range("B695").Value = 'Vega Monumental Concepción'
$ws.Range("C695").Value = 'Bíobío'
$ws.Range("D695").Value = 45041
$ws.Range("E695").Value = 8
$ws.Range("F695").Value = 100112020
$ws.Range("G695").Value = 'Tomate'
$ws.Range("H695").Value = 'Larga vida'
$ws.Range("I695").Value = 'Primera'
$ws.Range("J695").Value = 550
$ws.Range("K695").Value = 14000
$ws.Range("L695").Value = 15000
$ws.Range("M695").Value = 14545
$ws.Range("N695").Value = '$/bandeja 18 kilos'
$ws.Range("O695").Value = 'Provincia de Quillota'
$ws.Range("P695").Value = 808
$ws.Range("Q695").Value = 18
$ws.Range("R695").Value = 'Hortaliza'
